$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (pushes "SUPOLACK HAIR SHAMPOO 200 ML" and every
# product below it down by one row) to make room for a new product line:
# "SULBIN 750MG VIAL".
$ws.Rows(8).Insert()

# Mirror the formatting used by the other product rows onto the freshly
# inserted (blank) row 8: copy the style from the row above (row 7), then
# merge the label/value groups the same way every other data row is merged.
$ws.Range("A7:N7").Copy()
$ws.Range("A8:N8").PasteSpecial(-4122)

$ws.Range("B8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()

# Fill in the new product's data.
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "SULBIN 750MG VIAL"
$ws.Range("H8").Value = "6:0"
$ws.Range("L8").Value = 35
$ws.Range("N8").Value = "1:0"

# Column A is just the running row counter (1, 2, 3, ...) - restore the
# sequential numbering for every row beneath the new one, since the insert
# operation shifted the old literal numbers down along with everything else.
$ws.Range("A9").Value = 6
$ws.Range("A10").Value = 7
$ws.Range("A11").Value = 8
$ws.Range("A12").Value = 9
$ws.Range("A13").Value = 10
$ws.Range("A14").Value = 11
$ws.Range("A15").Value = 12
$ws.Range("A16").Value = 13

# Update the totals row (now on row 17) to include the new product's price.
$ws.Range("K17").Value = 692.5

# Row heights are fixed per row position in this template - put every row
# back to its usual height (the insert operation shifted the heights down
# together with the old values) and give the newly added rows their own
# auto-fitted heights.
$ws.Rows(8).RowHeight = 25.5
$ws.Rows(9).RowHeight = 24.75
$ws.Rows(10).RowHeight = 25.5
$ws.Rows(11).RowHeight = 24.75
$ws.Rows(12).RowHeight = 25.5
$ws.Rows(13).RowHeight = 25.5
$ws.Rows(14).RowHeight = 24.75
$ws.Rows(15).RowHeight = 25.5
$ws.Rows(16).RowHeight = 24.75
$ws.Rows(17).RowHeight = 26.25
$ws.Rows(18).RowHeight = 16.5
